$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 942.871
$ws.Range("J17").Value = 950.98334
$ws.Range("L17").Value = 2852.95002
$ws.Range("N17").Value = -3188.95002

$ws.Range("H58").Value = 2089.0417
$ws.Range("J58").Value = 3249.5715
$ws.Range("L58").Value = 9748.7145
$ws.Range("N58").Value = -10048.7145

$ws.Range("H61").Value = 5007.5
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30344

$ws.Range("H82").Value = 5125
$ws.Range("I82").Value = 250
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 750
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -344
$ws.Range("N82").Value = -30812

$ws.Range("H85").Value = 5125
$ws.Range("I85").Value = 250
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 750
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = 654
$ws.Range("N85").Value = -32808

$ws.Range("H99").Value = 365.36365
$ws.Range("I99").Value = 283.8
$ws.Range("J99").Value = 1181
$ws.Range("K99").Value = 851.4000000000001
$ws.Range("L99").Value = 3543
$ws.Range("M99").Value = 646.5999999999999
$ws.Range("N99").Value = -6539

$ws.Range("H101").Value = 403.66666
$ws.Range("J101").Value = 600
$ws.Range("L101").Value = 1800
$ws.Range("N101").Value = -5044

$ws.Range("H104").Value = 141.5
$ws.Range("I104").Value = 141.5
$ws.Range("K104").Value = 424.5
$ws.Range("M104").Value = 1322.5

$ws.Range("H115").Value = 415
$ws.Range("I115").Value = 415
$ws.Range("K115").Value = 1245
$ws.Range("M115").Value = 322

$ws.Range("H118").Value = 1029.25
$ws.Range("I118").Value = 640.5714
$ws.Range("K118").Value = 1921.7142
$ws.Range("M118").Value = -264.7142000000001

$ws.Range("H127").Value = 750.6667
$ws.Range("I127").Value = 377.25
$ws.Range("K127").Value = 1131.75
$ws.Range("M127").Value = 3828.25

$ws.Range("H129").Value = 1303.125
$ws.Range("I129").Value = 795
$ws.Range("J129").Value = 2150
$ws.Range("K129").Value = 2385
$ws.Range("L129").Value = 6450
$ws.Range("M129").Value = 2615
$ws.Range("N129").Value = -16450

$ws.Range("H132").Value = 2709693.8
$ws.Range("I132").Value = 2930832.5
$ws.Range("K132").Value = 8792497.5
$ws.Range("M132").Value = -8789967.5

$ws.Range("H137").Value = 10338.611
$ws.Range("I137").Value = 12416.5
$ws.Range("J137").Value = 3066
$ws.Range("K137").Value = 37249.5
$ws.Range("L137").Value = 9198
$ws.Range("M137").Value = -34699.5
$ws.Range("N137").Value = -14298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11990.637
$ws.Range("I61").Value = 1066
$ws.Range("J61").Value = 16087.375
$ws.Range("K61").Value = 1066
$ws.Range("L61").Value = 16087.375
$ws.Range("M61").Value = -854
$ws.Range("N61").Value = -16511.375

$ws.Range("H97").Value = 1504.8077
$ws.Range("I97").Value = 1101.875
$ws.Range("K97").Value = 1101.875
$ws.Range("M97").Value = -605.875

$ws.Range("H136").Value = 11990.637
$ws.Range("I136").Value = 1066
$ws.Range("J136").Value = 16087.375
$ws.Range("K136").Value = 3198
$ws.Range("L136").Value = 48262.125
$ws.Range("M136").Value = -648
$ws.Range("N136").Value = -53362.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1165.5
$ws.Range("I86").Value = 1135.2059
$ws.Range("J86").Value = 1251.3334
$ws.Range("K86").Value = 1135.2059
$ws.Range("L86").Value = 1251.3334
$ws.Range("M86").Value = -12.20589999999993
$ws.Range("N86").Value = -3497.3334

$ws.Range("H89").Value = 1165.5
$ws.Range("I89").Value = 1135.2059
$ws.Range("J89").Value = 1251.3334
$ws.Range("K89").Value = 5676.0295
$ws.Range("L89").Value = 6256.666999999999
$ws.Range("M89").Value = -60.02949999999964
$ws.Range("N89").Value = -17488.667

$ws.Range("H94").Value = 2042.0714
$ws.Range("I94").Value = 1579
$ws.Range("J94").Value = 3199.75
$ws.Range("K94").Value = 1579
$ws.Range("L94").Value = 3199.75
$ws.Range("M94").Value = -1128
$ws.Range("N94").Value = -4101.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 22625
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 60166.668
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 60166.668
$ws.Range("M4").Value = 12
$ws.Range("N4").Value = -60390.668

$ws.Range("H7").Value = 78.07143000000001
$ws.Range("I7").Value = 91.181816
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 91.181816
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = 21.818184
$ws.Range("N7").Value = -256

$ws.Range("H22").Value = 750.0909
$ws.Range("I22").Value = 600.5
$ws.Range("J22").Value = 835.5714
$ws.Range("K22").Value = 600.5
$ws.Range("L22").Value = 835.5714
$ws.Range("M22").Value = -250.5
$ws.Range("N22").Value = -1535.5714

$ws.Range("H55").Value = 26444.334
$ws.Range("I55").Value = 8500
$ws.Range("J55").Value = 62333
$ws.Range("K55").Value = 8500
$ws.Range("L55").Value = 62333
$ws.Range("M55").Value = -8185
$ws.Range("N55").Value = -62963

$ws.Range("H70").Value = 35300
$ws.Range("J70").Value = 35300
$ws.Range("L70").Value = 35300
$ws.Range("N70").Value = -35930

$ws.Range("H73").Value = 35300
$ws.Range("J73").Value = 35300
$ws.Range("L73").Value = 35300
$ws.Range("N73").Value = -37484

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8142
$ws.Range("I3").Value = 3399.4
$ws.Range("K3").Value = 10198.2
$ws.Range("M3").Value = -10086.2

$ws.Range("H4").Value = 83168500
$ws.Range("I4").Value = 58326764
$ws.Range("J4").Value = 154144880
$ws.Range("K4").Value = 174980292
$ws.Range("L4").Value = 462434640
$ws.Range("M4").Value = -174980180
$ws.Range("N4").Value = -462434864

$ws.Range("H20").Value = 1250
$ws.Range("I20").Value = 500
$ws.Range("K20").Value = 1500
$ws.Range("M20").Value = -1273

$ws.Range("H21").Value = 637.75
$ws.Range("J21").Value = 396.5
$ws.Range("L21").Value = 1189.5
$ws.Range("N21").Value = -1535.5

$ws.Range("H22").Value = 4435.643
$ws.Range("I22").Value = 399.33334
$ws.Range("K22").Value = 1198.00002
$ws.Range("M22").Value = -1029.00002

$ws.Range("H25").Value = 2590.5715
$ws.Range("I25").Value = 3182.75
$ws.Range("K25").Value = 9548.25
$ws.Range("M25").Value = -9379.25

$ws.Range("H27").Value = 4435.643
$ws.Range("I27").Value = 399.33334
$ws.Range("K27").Value = 1198.00002
$ws.Range("M27").Value = -1096.00002

$ws.Range("H30").Value = 2590.5715
$ws.Range("I30").Value = 3182.75
$ws.Range("K30").Value = 9548.25
$ws.Range("M30").Value = -9446.25

$ws.Range("H34").Value = 4016.7646
$ws.Range("J34").Value = 4998.846
$ws.Range("L34").Value = 14996.538
$ws.Range("N34").Value = -15164.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6935.4443
$ws.Range("I70").Value = 6682.0625
$ws.Range("K70").Value = 6682.0625
$ws.Range("M70").Value = -6412.0625

$ws.Range("H73").Value = 6935.4443
$ws.Range("I73").Value = 6682.0625
$ws.Range("K73").Value = 6682.0625
$ws.Range("M73").Value = -5746.0625

$ws.Range("H80").Value = 17083.334
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 17083.334
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 17083.334
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -19079.334

$ws.Range("H83").Value = 17083.334
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 17083.334
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 85416.67
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -95400.67

$ws.Range("H97").Value = 1489.3448
$ws.Range("I97").Value = 951.3889
$ws.Range("J97").Value = 2369.6365
$ws.Range("K97").Value = 951.3889
$ws.Range("L97").Value = 2369.6365
$ws.Range("M97").Value = -455.3889
$ws.Range("N97").Value = -3361.6365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1198.3846
$ws.Range("J22").Value = 1389
$ws.Range("L22").Value = 1389
$ws.Range("N22").Value = -1979

$ws.Range("H27").Value = 1198.3846
$ws.Range("J27").Value = 1389
$ws.Range("L27").Value = 1389
$ws.Range("N27").Value = -1603

$ws.Range("H40").Value = 1771.5588
$ws.Range("I40").Value = 1673.7273
$ws.Range("K40").Value = 1673.7273
$ws.Range("M40").Value = -1537.7273

$ws.Range("H132").Value = 1439.4286
$ws.Range("I132").Value = 901.3889
$ws.Range("J132").Value = 4667.6665
$ws.Range("K132").Value = 2704.1667
$ws.Range("L132").Value = 14002.9995
$ws.Range("M132").Value = -174.1667000000002
$ws.Range("N132").Value = -19062.9995

$ws.Range("H136").Value = 2917.0264
$ws.Range("I136").Value = 3082.2593
$ws.Range("J136").Value = 2511.4546
$ws.Range("K136").Value = 9246.777900000001
$ws.Range("L136").Value = 7534.3638
$ws.Range("M136").Value = -6696.777900000001
$ws.Range("N136").Value = -12634.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 120000
$ws.Range("J70").Value = 120000
$ws.Range("L70").Value = 120000
$ws.Range("N70").Value = -120630

$ws.Range("H73").Value = 120000
$ws.Range("J73").Value = 120000
$ws.Range("L73").Value = 120000
$ws.Range("N73").Value = -122184

$ws.Range("H96").Value = 1191
$ws.Range("J96").Value = 1515.5
$ws.Range("L96").Value = 1515.5
$ws.Range("N96").Value = -4261.5
